$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "cultype" column (C) values "monocul1", "monocul2", "cocul" are being split
# into two more specific labels each, depending on whether the row's "group"
# entry (column A) is the original "... COLLECTION" library stock row (first of
# each replicate pair) or the "remade" row (second of the pair, no COLLECTION
# suffix):
#   monocul1 -> monocul1 library / monocul1 remade
#   monocul2 -> monocul2 library / monocul2 remade
#   cocul    -> cocul library    / cocul remade
#
# Data rows are 4..75. Each pass below only touches one (old value, is
# COLLECTION row) combination at a time and scans top-to-bottom, which
# reproduces the original authoring order of the new shared-string table.

$lastRow = 75

# Pass 1: monocul1 -> "monocul1 library" (COLLECTION rows)
for ($r = 4; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "monocul1" -and $a -like "*COLLECTION*") {
        $ws.Cells.Item($r, 3).Value = "monocul1 library"
    }
}

# Pass 2: monocul2 -> "monocul2 library" (COLLECTION rows)
for ($r = 4; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "monocul2" -and $a -like "*COLLECTION*") {
        $ws.Cells.Item($r, 3).Value = "monocul2 library"
    }
}

# Pass 3: monocul1 -> "monocul1 remade" (non-COLLECTION rows)
for ($r = 4; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "monocul1" -and $a -notlike "*COLLECTION*") {
        $ws.Cells.Item($r, 3).Value = "monocul1 remade"
    }
}

# Pass 4: monocul2 -> "monocul2 remade" (non-COLLECTION rows)
for ($r = 4; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "monocul2" -and $a -notlike "*COLLECTION*") {
        $ws.Cells.Item($r, 3).Value = "monocul2 remade"
    }
}

# Pass 5: cocul -> "cocul library" (COLLECTION rows)
for ($r = 4; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "cocul" -and $a -like "*COLLECTION*") {
        $ws.Cells.Item($r, 3).Value = "cocul library"
    }
}

# Pass 6: cocul -> "cocul remade" (non-COLLECTION rows)
for ($r = 4; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    if ($c -eq "cocul" -and $a -notlike "*COLLECTION*") {
        $ws.Cells.Item($r, 3).Value = "cocul remade"
    }
}

# The longer labels in column C no longer fit the old column width, so widen it
# (mirrors the auto-adjusted column width seen after the edit).
$ws.Columns.Item(3).ColumnWidth = 21.33

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("C79").Select()
